# Implement change requests for forms
#
# 1. On the "survey" sheet, remove the "begin screen" / "end screen" marker
#    rows (clause rows with no type/name set) that wrapped the q34-q37
#    screen questions, so the question rows shift up and no longer sit
#    inside an explicit begin/end screen block.
# 2. Bump the form_version on the "settings" sheet to reflect the new
#    form revision.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

# Row 2 holds the "begin screen" clause (column A only) - delete it first.
$survey.Rows.Item(2).Delete()

# After the row-2 deletion, the old row 7 ("end screen") has shifted up to
# row 6 - delete that marker row too.
$survey.Rows.Item(6).Delete()

# Update the form_version setting to the new revision number.
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B3").Value = 20210304001
